# Station Standard.xlsx - data corrections on the "Station" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Station")
$ws.Activate()

# Row 20 (REFILL STATION B19): Drawers to count / KPI corrections
$ws.Range("D20").Value = 9
$ws.Range("F20").Value = 12.5

# Row 21 (REFILL STATION B20): KPI correction
$ws.Range("F21").Value = 16.5

# Row 26 (REFILL STATION C25) and Row 34 (REFILL STATION C33) had their
# Type/Config/metrics swapped between each other.
$ws.Range("B26").Value = "Atlas Box & Bond Bags"
$ws.Range("C26").Value = "Atlas Box & Bond Bags"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 0

$ws.Range("B34").Value = "Service Cart"
$ws.Range("C34").Value = "Service Cart 1"
$ws.Range("D34").Value = 5
$ws.Range("E34").Value = 14
$ws.Range("F34").Value = 9.6

# Restore the selection/scroll state recorded in the saved workbook view
$ws.Range("D35").Select()
